$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AB (2..28) hold the data that gets shuffled between rows; column A
# (the running match index) never moves.
$cols = 2..28

# Row permutation: for each destination row, the column letter-value it should
# end up holding is a full copy of another row's current (pre-edit) contents.
# These come from matching the unique match-id in column B across the diff.
$moves = [ordered]@{
    15  = 17
    16  = 15
    17  = 16
    170 = 171
    171 = 170
    175 = 177
    176 = 178
    177 = 176
    178 = 175
    237 = 238
    238 = 239
    239 = 237
}

# Snapshot every source row BEFORE any writes happen, since several rows are
# both a source and a destination in this permutation (cyclic shuffles).
$snapshots = @{}
foreach ($destRow in $moves.Keys) {
    $srcRow = $moves[$destRow]
    if (-not $snapshots.ContainsKey($srcRow)) {
        $rowData = @{}
        foreach ($c in $cols) {
            $rowData[$c] = $ws.Cells.Item($srcRow, $c).Value2
        }
        $snapshots[$srcRow] = $rowData
    }
}

# Now write the snapshotted source data into each destination row.
foreach ($destRow in $moves.Keys) {
    $srcRow = $moves[$destRow]
    $rowData = $snapshots[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $rowData[$c]
    }
}
